$wb = $excel.ActiveWorkbook

# --- Energy_prices: District_Heating connection_flow_cost column (D6:D53) ---
# Previously hard-coded to -1 for every timestep; fix to 0.
$wsPrices = $wb.Worksheets.Item("Energy_prices")
$wsPrices.Range("D6:D53").Value = 0

# --- Object__node_node: fix bug in connection out_in relation ---
# Add the missing Electrolyzer -> Hydrogen_Kasso/Waste_Heat unit__node__node
# relation (row 7), and add the missing pipeline_District_Heating
# connection__node__node relation (District_Heating -> Waste_Heat), which
# pushes the existing storage connection rows down.
$wsRel = $wb.Worksheets.Item("Object__node_node")

$wsRel.Rows.Item(7).Insert()
$wsRel.Range("A7").Value = "unit__node__node"
$wsRel.Range("B7").Value = "unit"
$wsRel.Range("C7").Value = "Electrolyzer"
$wsRel.Range("D7").Value = "Hydrogen_Kasso"
$wsRel.Range("E7").Value = "Waste_Heat"
$wsRel.Range("F7").Value = "fix_ratio_out_out_unit_flow"
$wsRel.Range("G7").Value = 1

$wsRel.Rows.Item(12).Insert()
$wsRel.Range("A12").Value = "connection__node__node"
$wsRel.Range("B12").Value = "connection"
$wsRel.Range("C12").Value = "pipeline_District_Heating"
$wsRel.Range("D12").Value = "District_Heating"
$wsRel.Range("E12").Value = "Waste_Heat"
$wsRel.Range("F12").Value = "fix_ratio_out_in_connection_flow"
$wsRel.Range("G12").Value = 1
